$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Miller", "SG,SF,PF", "Charlotte Hornets"),
    @("Duncan Robinson", "SG,SF", "Miami Heat"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Herbert Jones", "SF,PF", "New Orleans Pelicans"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Obi Toppin", "PF", "Indiana Pacers"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Ben Simmons", "PG,C", "Brooklyn Nets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Damian Lillard", "PG", "Milwaukee Bucks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
